$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.147234082221985
$ws.Range("B1").Value = 2.24649453163147
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.154542684555054
$ws.Range("E1").Value = 1.067817211151123
